# Apply updated crypto price/volume data (and restore row 5/6 coin order)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.755.20"
$ws.Range("E2").Value = "  -3.75%  "

$ws.Range("D3").Value = "1.616.26"
$ws.Range("E3").Value = "  -3.60%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'306.69"
$ws.Range("E6").Value = "  -2.23%  "

$ws.Range("D7").Value = "'0.3921"
$ws.Range("E7").Value = "  -0.51%  "

$ws.Range("D8").Value = "'0.3836"
$ws.Range("E8").Value = "  -2.84%  "

$ws.Range("D9").Value = "'1.004"
$ws.Range("E9").Value = "  +0.36%  "

$ws.Range("D10").Value = "'1.364"
$ws.Range("E10").Value = "  -2.93%  "

$ws.Range("D11").Value = "'49.74"
$ws.Range("E11").Value = "  -2.35%  "

$ws.Range("D12").Value = "'0.08431"
$ws.Range("E12").Value = "  -2.53%  "

$ws.Range("D13").Value = "'23.79"
$ws.Range("E13").Value = "  -5.82%  "

$ws.Range("D14").Value = "'7.036"
$ws.Range("E14").Value = "  -4.10%  "

$ws.Range("D15").Value = "'7.530"
$ws.Range("E15").Value = "  -2.46%  "

$ws.Range("D16").Value = "'0.00001275"
$ws.Range("E16").Value = "  -3.07%  "

$ws.Range("D17").Value = "1.617.62"
$ws.Range("E17").Value = "  -3.54%  "

$ws.Range("D18").Value = "'93.75"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("D19").Value = "'0.06923"
$ws.Range("E19").Value = "  -1.35%  "

$ws.Range("D20").Value = "'19.99"
$ws.Range("E20").Value = "  -5.48%  "

$ws.Range("D21").Value = "'6.805"
$ws.Range("E21").Value = "  -3.96%  "

$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").Value = "'13.37"
$ws.Range("E23").Value = "  -4.25%  "

$ws.Range("D24").Value = "23.766.92"
$ws.Range("E24").Value = "  -3.67%  "

$ws.Range("D25").Value = "'2.457"
$ws.Range("E25").Value = "  +4.61%  "

$ws.Range("D26").Value = "'2.852"
$ws.Range("E26").Value = "  +2.01%  "

$ws.Range("D27").Value = "'22.21"
$ws.Range("E27").Value = "  -3.69%  "

$ws.Range("D28").Value = "'156.22"
$ws.Range("E28").Value = "  -2.45%  "

$ws.Range("D29").Value = "'139.98"
$ws.Range("E29").Value = "  -4.30%  "

$ws.Range("D30").Value = "'5.265"
$ws.Range("E30").Value = "  -10.08%  "

$ws.Range("D31").Value = "'7.839"
$ws.Range("E31").Value = "  -5.92%  "

$ws.Range("D32").Value = "'2.505"
$ws.Range("E32").Value = "  +0.34%  "

$ws.Range("D33").Value = "1.796.95"
$ws.Range("E33").Value = "  -3.56%  "

$ws.Range("D34").Value = "'0.08123"
$ws.Range("E34").Value = "  -1.79%  "

$ws.Range("D35").Value = "'0.9793"
$ws.Range("E35").Value = "  -1.22%  "

$ws.Range("D36").Value = "'0.02874"
$ws.Range("E36").Value = "  -6.89%  "

$ws.Range("D37").Value = "'6.567"
$ws.Range("E37").Value = "  -5.98%  "

$ws.Range("D38").Value = "'0.2664"
$ws.Range("E38").Value = "  -5.35%  "

$ws.Range("D39").Value = "'0.09127"
$ws.Range("E39").Value = "  -5.36%  "

$ws.Range("D40").Value = "'10.26"
$ws.Range("E40").Value = "  -0.43%  "

$ws.Range("D41").Value = "'13.48"
$ws.Range("E41").Value = "  -0.28%  "

$ws.Range("D42").Value = "'1.422"
$ws.Range("E42").Value = "  -6.57%  "

$ws.Range("D43").Value = "'0.7484"
$ws.Range("E43").Value = "  -5.54%  "

$ws.Range("D44").Value = "'16.04"
$ws.Range("E44").Value = "  -3.59%  "

$ws.Range("D45").Value = "'0.6878"
$ws.Range("E45").Value = "  -3.29%  "

$ws.Range("D46").Value = "'2.461"
$ws.Range("E46").Value = "  -4.20%  "

$ws.Range("D47").Value = "'4.061"
$ws.Range("E47").Value = "  -2.70%  "

$ws.Range("D48").Value = "'1.000"
$ws.Range("E48").Value = "  +0.04%  "

$ws.Range("D49").Value = "'0.08232"
$ws.Range("E49").Value = "  -5.17%  "

$ws.Range("D50").Value = "'133.64"
$ws.Range("E50").Value = "  -3.21%  "

$ws.Range("D51").Value = "'1.212"
$ws.Range("E51").Value = "  -8.76%  "

# Reset styling on cells that needed a text-forcing quote prefix, so no
# spurious number formats / quote-prefix styling linger on the cells.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

